$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for the rows that changed.
$updates = @{
    "F4"  = 4
    "F5"  = -3
    "F6"  = -1
    "F7"  = 3
    "F13" = 0
    "F22" = -4
    "F24" = -2
    "F25" = 2
    "F26" = 4
    "F27" = 0
    "F29" = -1
    "F30" = 0
    "F43" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
